$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.734859704971313
$ws.Range("B1").Value = 2.665366649627686
$ws.Range("C1").Value = 2.190991163253784
$ws.Range("D1").Value = 2.07158899307251
$ws.Range("E1").Value = 1.798366189002991
